$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.20623640891663
$ws.Range("C2").Value = 6.385358286240186
$ws.Range("D2").Value = 6.031926204540781
$ws.Range("E2").Value = 11.14999939463377
$ws.Range("G2").Value = 51.99776718609861
$ws.Range("H2").Value = 19.92052412535984
$ws.Range("K2").Value = 13.4591484081665
$ws.Range("L2").Value = 10.0609563846518
$ws.Range("N2").Value = 22.39055922143823
$ws.Range("B3").Value = 17.00846135065884
$ws.Range("C3").Value = 6.291716114484077
$ws.Range("D3").Value = 5.921587917251717
$ws.Range("E3").Value = 11.15722288934395
$ws.Range("G3").Value = 51.81089791319044
$ws.Range("H3").Value = 19.93590165317431
$ws.Range("K3").Value = 13.32417713301473
$ws.Range("L3").Value = 10.05268370440948
$ws.Range("N3").Value = 22.43964158609521
$ws.Range("B4").Value = 16.89061543450823
$ws.Range("C4").Value = 6.232470128137936
$ws.Range("D4").Value = 5.854646153820521
$ws.Range("E4").Value = 11.16344942164739
$ws.Range("G4").Value = 51.70839836525851
$ws.Range("H4").Value = 19.94865658522837
$ws.Range("K4").Value = 13.24421737541546
$ws.Range("L4").Value = 10.04949064899785
$ws.Range("N4").Value = 22.47168246559633
$ws.Range("B5").Value = 16.84355112225058
$ws.Range("C5").Value = 6.20789772226456
$ws.Range("D5").Value = 5.827609801962024
$ws.Range("E5").Value = 11.16643740675645
$ws.Range("G5").Value = 51.66973005462512
$ws.Range("H5").Value = 19.95468667366077
$ws.Range("K5").Value = 13.21240139875742
$ws.Range("L5").Value = 10.04866510430183
$ws.Range("N5").Value = 22.48521846208448
$ws.Range("B6").Value = 16.83579565057697
$ws.Range("C6").Value = 6.203791891331231
$ws.Range("D6").Value = 5.823136283894241
$ws.Range("E6").Value = 11.16696077911949
$ws.Range("G6").Value = 51.66349711934465
$ws.Range("H6").Value = 19.95573821092913
$ws.Range("K6").Value = 13.20716579480936
$ws.Range("L6").Value = 10.04855678395487
$ws.Range("N6").Value = 22.48749504745281
$ws.Range("B7").Value = 16.88997675334056
$ws.Range("C7").Value = 6.232140459404481
$ws.Range("D7").Value = 5.854280495302461
$ws.Range("E7").Value = 11.16348789393915
$ws.Range("G7").Value = 51.70786428670618
$ws.Range("H7").Value = 19.94873454016594
$ws.Range("K7").Value = 13.24378513730724
$ws.Range("L7").Value = 10.04947758800866
$ws.Range("N7").Value = 22.47186307668101
$ws.Range("B8").Value = 17.13733316809688
$ws.Range("C8").Value = 6.353439822011945
$ws.Range("D8").Value = 5.993738290761089
$ws.Range("E8").Value = 11.15211842395117
$ws.Range("G8").Value = 51.93080972953873
$ws.Range("H8").Value = 19.9251382801977
$ws.Range("K8").Value = 13.41202760162009
$ws.Range("L8").Value = 10.05771330051778
$ws.Range("N8").Value = 22.40708778677241
$ws.Range("B9").Value = 17.64830685774235
$ws.Range("C9").Value = 6.576981099408499
$ws.Range("D9").Value = 6.27185698784542
$ws.Range("E9").Value = 11.1440239048796
$ws.Range("G9").Value = 52.4639100039028
$ws.Range("H9").Value = 19.90518057804282
$ws.Range("K9").Value = 13.76342390805795
$ws.Range("L9").Value = 10.08876352145035
$ws.Range("N9").Value = 22.29515929009064
$ws.Range("B10").Value = 18.03598803476261
$ws.Range("C10").Value = 6.731919661058328
$ws.Range("D10").Value = 6.476712003812188
$ws.Range("E10").Value = 11.14671263283031
$ws.Range("G10").Value = 52.91218753140951
$ws.Range("H10").Value = 19.9065864341752
$ws.Range("K10").Value = 14.03240752358355
$ws.Range("L10").Value = 10.12055316329339
$ws.Range("N10").Value = 22.22211041388391
$ws.Range("B11").Value = 18.21425687554213
$ws.Range("C11").Value = 6.800276163709622
$ws.Range("D11").Value = 6.569550649511208
$ws.Range("E11").Value = 11.14980378984019
$ws.Range("G11").Value = 53.12794204978386
$ws.Range("H11").Value = 19.91071614374776
$ws.Range("K11").Value = 14.15662348157089
$ws.Range("L11").Value = 10.13693590904992
$ws.Range("N11").Value = 22.19086960651929
$ws.Range("B12").Value = 18.28197383499878
$ws.Range("C12").Value = 6.825846573079801
$ws.Range("D12").Value = 6.604618703622503
$ws.Range("E12").Value = 11.15124207618783
$ws.Range("G12").Value = 53.21129886776218
$ws.Range("H12").Value = 19.91278140510717
$ws.Range("K12").Value = 14.20388515062595
$ws.Range("L12").Value = 10.14341302890991
$ws.Range("N12").Value = 22.17932544154035
$ws.Range("B13").Value = 18.26738143425678
$ws.Range("C13").Value = 6.820353663470116
$ws.Range("D13").Value = 6.597070658994024
$ws.Range("E13").Value = 11.15092042370454
$ws.Range("G13").Value = 53.19327363228077
$ws.Range("H13").Value = 19.91231432291119
$ws.Range("K13").Value = 14.19369723225975
$ws.Range("L13").Value = 10.14200595923014
$ws.Range("N13").Value = 22.18179896217649
$ws.Range("B14").Value = 18.21982413937968
$ws.Range("C14").Value = 6.802386213995486
$ws.Range("D14").Value = 6.572437689924642
$ws.Range("E14").Value = 11.14991675837092
$ws.Range("G14").Value = 53.13476696582183
$ws.Range("H14").Value = 19.9108760081803
$ws.Range("K14").Value = 14.16050748745416
$ws.Range("L14").Value = 10.13746332571361
$ws.Range("N14").Value = 22.18991412837049
$ws.Range("B15").Value = 18.19071942845439
$ws.Range("C15").Value = 6.791339387696094
$ws.Range("D15").Value = 6.557336752606978
$ws.Range("E15").Value = 11.14933682218987
$ws.Range("G15").Value = 53.0991440896593
$ws.Range("H15").Value = 19.91006028188036
$ws.Range("K15").Value = 14.14020568440943
$ws.Range("L15").Value = 10.13471633396439
$ws.Range("N15").Value = 22.19492215351684
$ws.Range("B16").Value = 18.02437068025228
$ws.Range("C16").Value = 6.727408843442078
$ws.Range("D16").Value = 6.470634476454343
$ws.Range("E16").Value = 11.14654813539526
$ws.Range("G16").Value = 52.89832167248024
$ws.Range("H16").Value = 19.9063867465246
$ws.Range("K16").Value = 14.02432336359548
$ws.Range("L16").Value = 10.11952092937072
$ws.Range("N16").Value = 22.22419210425785
$ws.Range("B17").Value = 17.9227645202661
$ws.Range("C17").Value = 6.687638477190948
$ws.Range("D17").Value = 6.417328359297874
$ws.Range("E17").Value = 11.14531522848835
$ws.Range("G17").Value = 52.77812206085362
$ws.Range("H17").Value = 19.90502683126005
$ws.Range("K17").Value = 13.95367768402329
$ws.Range("L17").Value = 10.11068914937618
$ws.Range("N17").Value = 22.24265779583974
$ws.Range("B18").Value = 17.8645070134451
$ws.Range("C18").Value = 6.664564030841852
$ws.Range("D18").Value = 6.386637197053675
$ws.Range("E18").Value = 11.14478198415573
$ws.Range("G18").Value = 52.71010295552766
$ws.Range("H18").Value = 19.90457319972588
$ws.Range("K18").Value = 13.91322122883078
$ws.Range("L18").Value = 10.10579039795472
$ws.Range("N18").Value = 22.25346603441291
$ws.Range("B19").Value = 17.84481553892383
$ws.Range("C19").Value = 6.656717423918815
$ws.Range("D19").Value = 6.376241560522784
$ws.Range("E19").Value = 11.14463167103905
$ws.Range("G19").Value = 52.68726597063677
$ws.Range("H19").Value = 19.90447604720572
$ws.Range("K19").Value = 13.89955508575353
$ws.Range("L19").Value = 10.10416294731089
$ws.Range("N19").Value = 22.2571576820513
$ws.Range("B20").Value = 17.93356212861629
$ws.Range("C20").Value = 6.691892823634331
$ws.Range("D20").Value = 6.423006360749945
$ws.Range("E20").Value = 11.14542827546288
$ws.Range("G20").Value = 52.79080228096014
$ws.Range("H20").Value = 19.90513759466587
$ws.Range("K20").Value = 13.96118002138057
$ws.Range("L20").Value = 10.11161058960655
$ws.Range("N20").Value = 22.24067271089895
$ws.Range("B21").Value = 18.23378766811752
$ws.Range("C21").Value = 6.807672300227152
$ws.Range("D21").Value = 6.579675666178191
$ws.Range("E21").Value = 11.15020430082685
$ws.Range("G21").Value = 53.15190725042036
$ws.Range("H21").Value = 19.9112848723394
$ws.Range("K21").Value = 14.17025038732571
$ws.Range("L21").Value = 10.13879021431909
$ws.Range("N21").Value = 22.18752274419074
$ws.Range("B22").Value = 18.43119964978916
$ws.Range("C22").Value = 6.88150263901247
$ws.Range("D22").Value = 6.681539509634791
$ws.Range("E22").Value = 11.15488574502086
$ws.Range("G22").Value = 53.39753491029639
$ws.Range("H22").Value = 19.91822479365681
$ws.Range("K22").Value = 14.30817460313785
$ws.Range("L22").Value = 10.15814526696261
$ws.Range("N22").Value = 22.15445348917711
$ws.Range("B23").Value = 18.32574866062889
$ws.Range("C23").Value = 6.842269053071511
$ws.Range("D23").Value = 6.627233141400513
$ws.Range("E23").Value = 11.15224475743893
$ws.Range("G23").Value = 53.26557412623229
$ws.Range("H23").Value = 19.91425365293792
$ws.Range("K23").Value = 14.23445841157303
$ws.Range("L23").Value = 10.14767053543486
$ws.Range("N23").Value = 22.1719506195635
$ws.Range("B24").Value = 17.92868003518929
$ws.Range("C24").Value = 6.689970086873521
$ws.Range("D24").Value = 6.420439472894957
$ws.Range("E24").Value = 11.14537661998425
$ws.Range("G24").Value = 52.78506616932773
$ws.Range("H24").Value = 19.90508649607829
$ws.Range("K24").Value = 13.9577877170779
$ws.Range("L24").Value = 10.11119344945017
$ws.Range("N24").Value = 22.24156956944048
$ws.Range("B25").Value = 17.50766216450002
$ws.Range("C25").Value = 6.518100889952996
$ws.Range("D25").Value = 6.196363512885069
$ws.Range("E25").Value = 11.14469506569454
$ws.Range("G25").Value = 52.30960210206231
$ws.Range("H25").Value = 19.90775878397354
$ws.Range("K25").Value = 13.66629482225015
$ws.Range("L25").Value = 10.07877830919009
$ws.Range("N25").Value = 22.32382472819644
